$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("2025-11-15", "2025-11-16", "2025-11-17", "2025-11-18", "2025-11-19", "2025-11-20", "2025-11-21")
$values = @(94092.49000000001, 94109.49000000001, 94704.02, 94864.87, 95273.05, 95748.95, 95687.57000000001)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
